$d = $word.ActiveDocument

# Locate the paragraph that currently reads " `t, `tde `tde `t."
# (the underline-tab date/place line) and the one right after it that
# reads "Local`tData".
$dateTarget = " `t, `tde `tde `t."
$localDataTarget = "Local`tData"

$dateParaIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq $dateTarget) {
        $dateParaIndex = $i
    }
    $i = $i + 1
}

if ($dateParaIndex -ge 0) {
    $datePara = $d.Paragraphs.Item($dateParaIndex + 1)

    # Add the new tab stop (4125 twips = 206.25 pt) to this paragraph.
    $datePara.Range.ParagraphFormat.TabStops.Add(206.25)

    # Replace the paragraph's run content (excluding the paragraph mark)
    # with the new text.
    $rng = $datePara.Range
    $rng.MoveEnd(1, -1)
    $rng.Text = "Colares " + [char]8211 + " PA, <DataCompleta>."

    # Tag the "DataCompleta" word with pt-BR language, matching the
    # auto-detected language Word would apply to the new text.
    $findRange = $datePara.Range.Duplicate
    $findRange.Find.ClearFormatting()
    $findRange.Find.Execute("DataCompleta")
    if ($findRange.Find.Found) {
        $findRange.LanguageID = "pt-BR"
    }
}

# Delete the "Local`tData" paragraph entirely (including its paragraph mark).
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq $localDataTarget) {
        $p.Range.Delete()
        break
    }
}
